$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "nemad" column (L) was populated with the full company name
# ("نفت سپاهان") for every data row; correct it to the stock ticker
# symbol ("شسپا") for rows 2 through 45.
$ws.Range("L2:L45").Value = "شسپا"

# Column L is now a touch narrower text ("شسپا" vs "نفت سپاهان"), so the
# author resized/best-fit column L (stored width 10) and left the
# selection on L8 when the file was saved.
$ws.Columns("L").ColumnWidth = 9.17
$ws.Range("L8").Select()
